$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold plain text values (e.g. "35.530.50",
# "  +3.09%  ") in the source data, not numbers/percentages. Force the text
# number format first so Excel does not auto-convert numeric-looking strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '35.530.50'
$ws.Cells.Item(3, 4).Value = '1.908.17'
$ws.Cells.Item(3, 5).Value = '  +3.09%  '
$ws.Cells.Item(4, 4).Value = '1.02'
$ws.Cells.Item(4, 5).Value = '  +0.56%  '
$ws.Cells.Item(5, 4).Value = '245.64'
$ws.Cells.Item(5, 5).Value = '  +3.28%  '
$ws.Cells.Item(6, 4).Value = '0.656'
$ws.Cells.Item(6, 5).Value = '  +5.57%  '
$ws.Cells.Item(7, 5).Value = '  +0.47%  '
$ws.Cells.Item(8, 4).Value = '41.83'
$ws.Cells.Item(8, 5).Value = '  -0.76%  '
$ws.Cells.Item(9, 5).Value = '  +5.36%  '
$ws.Cells.Item(10, 4).Value = '50.18'
$ws.Cells.Item(10, 5).Value = '  +8.01%  '
$ws.Cells.Item(11, 4).Value = '0.0716'
$ws.Cells.Item(11, 5).Value = '  +3.37%  '
$ws.Cells.Item(12, 5).Value = '  +1.06%  '
$ws.Cells.Item(13, 4).Value = '2.186.84'
$ws.Cells.Item(13, 5).Value = '  +3.15%  '
$ws.Cells.Item(14, 4).Value = '12.19'
$ws.Cells.Item(14, 5).Value = '  +7.01%  '
$ws.Cells.Item(15, 4).Value = '0.698'
$ws.Cells.Item(15, 5).Value = '  +3.72%  '
$ws.Cells.Item(16, 4).Value = '4.88'
$ws.Cells.Item(16, 5).Value = '  +2.52%  '
$ws.Cells.Item(17, 4).Value = '1.896.76'
$ws.Cells.Item(17, 5).Value = '  +2.58%  '
$ws.Cells.Item(18, 4).Value = '35.545.55'
$ws.Cells.Item(18, 5).Value = '  +1.47%  '
$ws.Cells.Item(19, 4).Value = '72.36'
$ws.Cells.Item(19, 5).Value = '  +3.21%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0819'
$ws.Cells.Item(20, 5).Value = '  +3.26%  '
$ws.Cells.Item(21, 4).Value = '245.03'
$ws.Cells.Item(21, 5).Value = '  +2.12%  '
$ws.Cells.Item(22, 4).Value = '12.61'
$ws.Cells.Item(22, 5).Value = '  +3.87%  '
$ws.Cells.Item(23, 4).Value = '4.80'
$ws.Cells.Item(23, 5).Value = '  +1.00%  '
$ws.Cells.Item(24, 5).Value = '  +0.47%  '
$ws.Cells.Item(25, 4).Value = '2.30'
$ws.Cells.Item(25, 5).Value = '  +0.97%  '
$ws.Cells.Item(26, 5).Value = '  +26.00%  '
$ws.Cells.Item(27, 4).Value = '170.99'
$ws.Cells.Item(27, 5).Value = '  +0.72%  '
$ws.Cells.Item(28, 4).Value = '8.40'
$ws.Cells.Item(28, 5).Value = '  +5.29%  '
$ws.Cells.Item(29, 4).Value = '18.31'
$ws.Cells.Item(29, 5).Value = '  +4.13%  '
$ws.Cells.Item(30, 4).Value = '0.127'
$ws.Cells.Item(30, 5).Value = '  +2.33%  '
$ws.Cells.Item(31, 4).Value = '4.16'
$ws.Cells.Item(31, 5).Value = '  +3.96%  '
$ws.Cells.Item(32, 4).Value = '0.0569'
$ws.Cells.Item(32, 5).Value = '  +2.44%  '
$ws.Cells.Item(33, 4).Value = '1.02'
$ws.Cells.Item(33, 5).Value = '  +0.55%  '
$ws.Cells.Item(34, 2).Value = 'ImmutableX'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(34, 4).Value = '0.931'
$ws.Cells.Item(34, 5).Value = '  +19.94%  '
$ws.Cells.Item(35, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(35, 4).Value = '4.15'
$ws.Cells.Item(35, 5).Value = '  +3.51%  '
$ws.Cells.Item(36, 4).Value = '1.74'
$ws.Cells.Item(36, 5).Value = '  +4.68%  '
$ws.Cells.Item(37, 4).Value = '2.05'
$ws.Cells.Item(37, 5).Value = '  +3.04%  '
$ws.Cells.Item(38, 4).Value = '1.33'
$ws.Cells.Item(38, 5).Value = '  +1.45%  '
$ws.Cells.Item(39, 4).Value = '0.0211'
$ws.Cells.Item(39, 5).Value = '  +4.78%  '
$ws.Cells.Item(40, 5).Value = '  +2.73%  '
$ws.Cells.Item(41, 4).Value = '0.0635'
$ws.Cells.Item(41, 5).Value = '  +14.19%  '
$ws.Cells.Item(42, 4).Value = '91.00'
$ws.Cells.Item(42, 5).Value = '  +1.02%  '
$ws.Cells.Item(43, 4).Value = '15.81'
$ws.Cells.Item(43, 5).Value = '  +7.59%  '
$ws.Cells.Item(44, 4).Value = '1.350.81'
$ws.Cells.Item(44, 5).Value = '  +0.29%  '
$ws.Cells.Item(45, 4).Value = '2.38'
$ws.Cells.Item(45, 5).Value = '  +2.37%  '
$ws.Cells.Item(46, 4).Value = '47.08'
$ws.Cells.Item(46, 5).Value = '  +36.63%  '
$ws.Cells.Item(47, 4).Value = '12.65'
$ws.Cells.Item(47, 5).Value = '  +3.82%  '
$ws.Cells.Item(48, 5).Value = '  +2.14%  '
$ws.Cells.Item(49, 5).Value = '  -0.15%  '
$ws.Cells.Item(50, 4).Value = '6.54'
$ws.Cells.Item(50, 5).Value = '  +0.17%  '
$ws.Cells.Item(51, 4).Value = '2.096.18'
$ws.Cells.Item(51, 5).Value = '  +2.86%  '
